$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "time_taken", matching the style of the existing header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Add time_taken values for each data row (F2:F7)
$ws.Range("F2").Value = "2021-10-05 13:39:35.764028"
$ws.Range("F3").Value = "2021-10-05 13:39:35.764039"
$ws.Range("F4").Value = "2021-10-05 13:39:35.764044"
$ws.Range("F5").Value = "2021-10-05 13:39:35.764047"
$ws.Range("F6").Value = "2021-10-05 13:39:35.764050"
$ws.Range("F7").Value = "2021-10-05 13:39:35.764053"
